$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Sebastián"
$ws.Range("B4").Value = "Palacio"
$ws.Range("C4").Value = 1000762620
$ws.Range("D4").Value = "sebasx200"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1234"
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value = "sebastian_palacio23231@elpoli,edu,co"
$ws.Range("G4").Value = "No tiene"
